$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsIea   = $wb.Worksheets.Item("IEA Data")
$wsHpebp = $wb.Worksheets.Item("HPEbP")

# --- "About" sheet: add state label + a "last updated" style date ---
$wsAbout.Range("B1").Value = "Oregon"
$wsAbout.Range("C1").Value = 44811
$wsAbout.Range("C1").NumberFormat = "mm-dd-yy"

# --- "HPEbP" sheet: electrolysis efficiency formula no longer divides by the extra 46 term ---
$wsHpebp.Range("B3").Formula = "=118/(162+2)"

# Normalize the now-redundant per-cell border flag on the "long term" block (Q:AI)
# so its style matches the plain numeric style used elsewhere on the sheet.
$wsHpebp.Range("Q2:AI6").Borders.LineStyle = -4142

# --- Selections / active sheet to match the saved UI state ---
$wsAbout.Range("B14").Select()
$wsIea.Range("D7:F7").Select()
$wsHpebp.Activate()
$wsHpebp.Range("C3").Select()
